$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9223046214701632
$ws.Range("J2").Value = 0.9223046214701632
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 23.896359208001
$ws.Range("R2").Value = 215.067232872009
$ws.Range("S2").Value = 0.2720401283092966
$ws.Range("T2").Value = 0.2720401283092966
$ws.Range("I3").Value = 0.9223046214701632
$ws.Range("J3").Value = 0.9223046214701632
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 45.19050474404901
$ws.Range("R3").Value = 406.7145426964411
$ws.Range("S3").Value = 0.5144562233068879
$ws.Range("T3").Value = 0.5144562233068879
$ws.Range("I4").Value = 0.9223046214701632
$ws.Range("J4").Value = 0.9223046214701632
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 11.92957531676367
$ws.Range("R4").Value = 107.366177850873
$ws.Range("S4").Value = 0.1358082698539786
$ws.Range("T4").Value = 0.1358082698539786
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.140557
$ws.Range("H5").Value = 0.421671
$ws.Range("I5").Value = 0.07769537852983674
$ws.Range("J5").Value = 0.07769537852983674
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 2.013040627717
$ws.Range("R5").Value = 18.117365649453
$ws.Range("S5").Value = 0.02291678936900992
$ws.Range("T5").Value = 0.02291678936900992
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.140557
$ws.Range("H6").Value = 0.421671
$ws.Range("I6").Value = 0.07769537852983674
$ws.Range("J6").Value = 0.07769537852983674
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("Q6").Value = 3.806869542133001
$ws.Range("R6").Value = 34.261825879197
$ws.Range("S6").Value = 0.04333803612861105
$ws.Range("T6").Value = 0.04333803612861105
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.140557
$ws.Range("H7").Value = 0.421671
$ws.Range("I7").Value = 0.07769537852983674
$ws.Range("J7").Value = 0.07769537852983674
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 1.004953079882333
$ws.Range("R7").Value = 9.044577718940999
$ws.Range("S7").Value = 0.01144055303221576
$ws.Range("T7").Value = 0.01144055303221576
